$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.435.43'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '3.974.50'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'616.67"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +11.88%  '
$ws.Range('D6').Value = "'169.23"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +11.09%  '
$ws.Range('D7').Value = "'0.681"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.29%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = "'0.760"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Value = "'0.187"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.55%  '
$ws.Range('D11').Value = "'55.54"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('D12').Value = "'0.0000338"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.51%  '
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').Value = '4.607.01'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').Value = '3.974.62'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').Value = "'14.07"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.87%  '
$ws.Range('D18').Value = "'20.42"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('D19').Value = '73.258.05'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = "'438.66"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.38%  '
$ws.Range('E22').Value = '  +12.68%  '
$ws.Range('D23').Value = "'95.96"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.22%  '
$ws.Range('E24').Value = '  -5.17%  '
$ws.Range('D25').Value = "'14.19"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('E26').Value = '  -3.70%  '
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('D28').Value = "'5.95"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = "'10.54"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.47%  '
$ws.Range('D30').Value = "'36.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.34%  '
$ws.Range('D31').Value = "'7.76"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('D32').Value = "'13.76"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('D33').Value = "'0.0000106"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +20.30%  '
$ws.Range('E34').Value = '  -4.13%  '
$ws.Range('D35').Value = "'48.07"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('D36').Value = "'71.09"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.52%  '
$ws.Range('D37').Value = "'648.40"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.67%  '
$ws.Range('D38').Value = "'0.431"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.81%  '
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('E41').Value = '  -2.44%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = "'3.23"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.94%  '
$ws.Range('D44').Value = "'0.0484"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('E45').Value = '  -5.78%  '
$ws.Range('D46').Value = "'3.20"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +42.37%  '
$ws.Range('D47').Value = "'0.149"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('E48').Value = '  +8.02%  '
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('E50').Value = '  -4.99%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = "'3.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.50%  '
